$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Goal: turn " kiMeSyic?" into " kiMeSy" + "n" + "?" as THREE separate
# runs (same bold rPr as before), while leaving the preceding run
# ("KMS KindMedSync KiMeSy") completely untouched, matching the target
# diff exactly.
#
# This COM-interop host re-coalesces same-formatted runs in a
# paragraph whenever a *text* edit touches that paragraph, but a pure
# *formatting* edit (e.g. toggling Bold) only splits the run(s) that
# the edit's range spans, without touching the rest of the paragraph.
# We exploit that: briefly detune the preceding run's Bold so it no
# longer matches the run we are editing (so the text edit cannot pull
# it into the merge), fix the text, restore Bold, then use a
# Bold-off/Bold-on no-op toggle on the single "n" character to carve
# it into its own run, separate from both " kiMeSy" and "?".
# ------------------------------------------------------------------

# 1) Locate the run right before the text we need to change.
$prefix = $d.Content
$null = $prefix.Find.Execute("KMS KindMedSync KiMeSy", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$prefixStart = $prefix.Start
$prefixEnd = $prefix.End

# 2) Detune its Bold so it no longer shares identical formatting with
#    the run we are about to edit (prevents an unwanted merge).
$prefixRange = $d.Range($prefixStart, $prefixEnd)
$prefixRange.Bold = 0

# 3) Find the exact "ic" inside " kiMeSyic?" that must become "n".
$target = $d.Content
$null = $target.Find.Execute("ic?", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$icStart = $target.Start
$icEnd = $target.End - 1   # exclude the trailing "?"

$icRange = $d.Range($icStart, $icEnd)
$icRange.Text = "n"

# 4) Restore the prefix run's Bold (back to its original look).
$prefixRange2 = $d.Range($prefixStart, $prefixEnd)
$prefixRange2.Bold = 1

# 5) Split the lone "n" away from " kiMeSy" and from "?" by toggling
#    Bold off then back on (a formatting no-op that nonetheless forces
#    a run boundary at both ends of the range).
$nFind = $d.Content
$null = $nFind.Find.Execute("n?", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$nStart = $nFind.Start
$nRange = $d.Range($nStart, $nStart + 1)
$nRange.Bold = 0
$nRange.Bold = 1

Write-Output "done"
